# Add two new market sheets (Netherlands, Denmark) to the workbook,
# copied from the "Italy" sheet template, and make "Denmark" the active tab.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Italy")

# --- Netherlands --------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet)
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2199"
$nl.Range("B2").Value = "Netherlands Market"
$nl.Range("D19").Select() | Out-Null

# --- Denmark -------------------------------------------------------------
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet2)
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"
$dk.Range("B4").Value = "NGC-2913/T2798"
$dk.Range("B2").Value = "Denmark Market"
$dk.Range("D19").Select() | Out-Null
